# EDIT INDEX & BURNDOWN CHART
# Reproduces:
#  - shift of the "Start" (C) column values down by one row for the task
#    table (rows 7-27 each take the C-value that used to belong to the row
#    above), leaving the D (task-name) column untouched per row
#  - a brand-new "Task 22.1" row inserted at row 28 (C28=2013, E28=3)
#  - the previous "Ideal"/"Actual" summary rows pushed down from 28/29 to
#    29/30, with their SUM ranges widened to include the new row 28
#  - a couple of independent effort-estimate edits (E6, P6, E23, E25)
#  - the burndown chart series ranges moving from row 28/29 to 29/30
#  - the selection/top-left view cell, and dropping the now-orphaned
#    B28:D28 / B29:D29 merges

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Push the "Actual" row (29 -> 30) and "Ideal" row (28 -> 29) down by
#    copying whole rows (this carries the correct cell styles with it),
#    bottom-up so we don't clobber source data before it is copied.
# ---------------------------------------------------------------------
$ws.Rows.Item(29).Copy($ws.Rows.Item(30))
$ws.Rows.Item(28).Copy($ws.Rows.Item(29))

# New task row 28 takes on the normal task-row formatting (copy a
# template row, e.g. row 27, then we overwrite its values below).
$ws.Rows.Item(27).Copy($ws.Rows.Item(28))

# ---------------------------------------------------------------------
# 2) Shift the "Start" column (C) down by one for every task row: each
#    row receives the C value that used to sit one row above it.
# ---------------------------------------------------------------------
$oldC = @{}
for ($r = 6; $r -le 27; $r++) {
    $oldC[$r] = $ws.Cells.Item($r, 3).Value()
}

$ws.Range("C7").Value = "Intro"
for ($r = 8; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = $oldC[$r - 1]
}

# ---------------------------------------------------------------------
# 3) New "Task 22.1" row content (row 28).
# ---------------------------------------------------------------------
$ws.Range("B28").Value = "History of Intel Processor"
$ws.Range("D28").Value = "Task 22.1"
$ws.Range("E28").Value = 3
$ws.Range("F28:Q28").ClearContents()

# ---------------------------------------------------------------------
# 4) Independent effort-estimate tweaks.
# ---------------------------------------------------------------------
$ws.Range("E6").Value = 4
$ws.Range("P6").Value = 0
$ws.Range("E23").Value = 4
$ws.Range("E25").Value = 3

# ---------------------------------------------------------------------
# 5) "Ideal" row, now at 29: E29-$E$29/12 chain instead of the old
#    SUM(E6:E27) total, and the totals row sums now reach down to the
#    new row 28.
# ---------------------------------------------------------------------
$ws.Range("E29").Formula = "=SUM(E6:E28)"
$ws.Range("F29").Formula = "=E29-`$E`$29/12"
$ws.Range("G29").Formula = "=F29-`$E`$29/12"
$ws.Range("H29").Formula = "=G29-`$E`$29/12"
$ws.Range("I29").Formula = "=H29-`$E`$29/12"
$ws.Range("J29").Formula = "=I29-`$E`$29/12"
$ws.Range("K29").Formula = "=J29-`$E`$29/12"
$ws.Range("L29").Formula = "=K29-`$E`$29/12"
$ws.Range("M29").Formula = "=L29-`$E`$29/12"
$ws.Range("N29").Formula = "=M29-`$E`$29/12"
$ws.Range("O29").Formula = "=N29-`$E`$29/12"
$ws.Range("P29").Formula = "=O29-`$E`$29/12"
$ws.Range("Q29").Formula = "=P29-`$E`$29/12"

# ---------------------------------------------------------------------
# 6) "Actual" row, now at 30: SUM(*6:*28) instead of SUM(*6:*27).
# ---------------------------------------------------------------------
$ws.Range("E30").Formula = "=SUM(E6:E28)"
$ws.Range("F30").Formula = "=SUM(F6:F28)"
$ws.Range("G30").Formula = "=SUM(G6:G28)"
$ws.Range("H30").Formula = "=SUM(H6:H28)"
$ws.Range("I30").Formula = "=SUM(I6:I28)"
$ws.Range("J30").Formula = "=SUM(J6:J28)"
$ws.Range("K30").Formula = "=SUM(K6:K28)"
$ws.Range("L30").Formula = "=SUM(L6:L28)"
$ws.Range("M30").Formula = "=SUM(M6:M28)"
$ws.Range("N30").Formula = "=SUM(N6:N28)"
$ws.Range("O30").Formula = "=SUM(O6:O28)"
$ws.Range("P30").Formula = "=SUM(P6:P28)"
$ws.Range("Q30").Formula = "=SUM(Q6:Q28)"

# ---------------------------------------------------------------------
# 7) The copy operations above dragged the B28:D28 / B29:D29 merges
#    along with them; the final layout has no merge on rows 28-30.
# ---------------------------------------------------------------------
$ws.Range("B28:D28").UnMerge()
$ws.Range("B29:D29").UnMerge()
$ws.Range("B30:D30").UnMerge()

# ---------------------------------------------------------------------
# 8) Burndown chart: series now read from rows 29/30 instead of 28/29.
# ---------------------------------------------------------------------
$chart = $ws.ChartObjects().Item(1).Chart
$chart.SeriesCollection().Item(1).Formula = '=SERIES("Ideal burndown",Sheet1!$F$5:$Q$5,Sheet1!$F$29:$Q$29,1)'
$chart.SeriesCollection().Item(2).Formula = '=SERIES("Actual burndown",Sheet1!$F$5:$Q$5,Sheet1!$F$30:$Q$30,2)'

# ---------------------------------------------------------------------
# 9) Selection / view.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("Q6").Select()
